# Update "想去人数" (want-to-go count, column F) figures across the
# workbook's sheets, and refresh the merged "全部类型" sheet (4th sheet)
# by removing duplicate rows that had been pulled in from the source
# sheets, then renumbering its serial-number column (A) and refreshing
# its F-column counts as well.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------
# Sheet "展览" (exhibitions): rows 2-16, column F = want-to-go counts
# ---------------------------------------------------------------
$ws1.Cells.Item(2,6).Value2  = 4454
$ws1.Cells.Item(3,6).Value2  = 2474
$ws1.Cells.Item(4,6).Value2  = 483
$ws1.Cells.Item(5,6).Value2  = 27
$ws1.Cells.Item(6,6).Value2  = 57
$ws1.Cells.Item(7,6).Value2  = 60
$ws1.Cells.Item(8,6).Value2  = 223
$ws1.Cells.Item(9,6).Value2  = 134
$ws1.Cells.Item(10,6).Value2 = 164
$ws1.Cells.Item(11,6).Value2 = 167
$ws1.Cells.Item(12,6).Value2 = 1670
$ws1.Cells.Item(13,6).Value2 = 302
$ws1.Cells.Item(14,6).Value2 = 3610
$ws1.Cells.Item(15,6).Value2 = 10
$ws1.Cells.Item(16,6).Value2 = 242

# ---------------------------------------------------------------
# Sheet "演出" (performances): rows 2-5, column F = want-to-go counts
# ---------------------------------------------------------------
$ws2.Cells.Item(2,6).Value2 = 1
$ws2.Cells.Item(3,6).Value2 = 44
$ws2.Cells.Item(4,6).Value2 = 3
$ws2.Cells.Item(5,6).Value2 = 7

# ---------------------------------------------------------------
# Sheet "全部类型" (all types, merges the above sheets): it currently
# carries 3 duplicated rows (old rows 5, 10 and 22) that need to be
# dropped. Delete bottom-to-top so earlier row numbers stay valid.
# ---------------------------------------------------------------
$ws4.Rows.Item(22).Delete()
$ws4.Rows.Item(10).Delete()
$ws4.Rows.Item(5).Delete()

# Refresh the want-to-go counts (col F) and the serial number (col A)
# for every remaining data row (now rows 2-20).
$ws4.Cells.Item(2,6).Value2   = 4454
$ws4.Cells.Item(2,1).Value2   = 1

$ws4.Cells.Item(3,6).Value2   = 2474
$ws4.Cells.Item(3,1).Value2   = 2

$ws4.Cells.Item(4,6).Value2   = 483
$ws4.Cells.Item(4,1).Value2   = 3

$ws4.Cells.Item(5,6).Value2   = 27
$ws4.Cells.Item(5,1).Value2   = 4

$ws4.Cells.Item(6,6).Value2   = 1
$ws4.Cells.Item(6,1).Value2   = 5

$ws4.Cells.Item(7,6).Value2   = 57
$ws4.Cells.Item(7,1).Value2   = 6

$ws4.Cells.Item(8,6).Value2   = 60
$ws4.Cells.Item(8,1).Value2   = 7

$ws4.Cells.Item(9,6).Value2   = 44
$ws4.Cells.Item(9,1).Value2   = 8

$ws4.Cells.Item(10,6).Value2  = 223
$ws4.Cells.Item(10,1).Value2  = 9

$ws4.Cells.Item(11,6).Value2  = 134
$ws4.Cells.Item(11,1).Value2  = 10

$ws4.Cells.Item(12,6).Value2  = 164
$ws4.Cells.Item(12,1).Value2  = 11

$ws4.Cells.Item(13,6).Value2  = 167
$ws4.Cells.Item(13,1).Value2  = 12

$ws4.Cells.Item(14,6).Value2  = 3
$ws4.Cells.Item(14,1).Value2  = 13

$ws4.Cells.Item(15,6).Value2  = 7
$ws4.Cells.Item(15,1).Value2  = 14

$ws4.Cells.Item(16,6).Value2  = 1670
$ws4.Cells.Item(16,1).Value2  = 15

$ws4.Cells.Item(17,6).Value2  = 302
$ws4.Cells.Item(17,1).Value2  = 16

$ws4.Cells.Item(18,6).Value2  = 3610
$ws4.Cells.Item(18,1).Value2  = 17

$ws4.Cells.Item(19,6).Value2  = 10
$ws4.Cells.Item(19,1).Value2  = 18

$ws4.Cells.Item(20,6).Value2  = 242
$ws4.Cells.Item(20,1).Value2  = 19
